$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the explicit (default-valued) left-to-right sheet view setting from
# the source workbook; harmless if the host doesn't surface it in the OOXML.
try { $excel.ActiveWindow.DisplayRightToLeft = $false } catch {}
try { $ws.DisplayRightToLeft = $false } catch {}

# Fix the existing row 2 "Problem" text.
$ws.Cells.Item(2, 6).Value = "Não consigo acessar minha conta"

# Append four more identical rows (3-6) with the same ticket data.
$rowValues = @("23/12/2024", "Thalles Gabriel", "Aluno", "Cajazeiras", "Polivalente", "Não consigo acessar minha conta")

for ($r = 3; $r -le 6; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
